# Update the "想去人数" (F column) figures on both the "展览" and
# "全部类型" worksheets to reflect the latest generated data.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 1585
    3  = 58
    6  = 71
    7  = 2775
    9  = 1786
    10 = 187
    11 = 82
    12 = 632
    15 = 161
    16 = 84
    17 = 89
    18 = 23
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Cells.Item($row, 6).Value = $updates[$row]
    }
}
